$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) - force text formatting so numeric-looking
# strings (e.g. thousand-separated prices) are preserved exactly as text
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '49.493.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.635.15'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '112.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '324.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.544'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.59'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.85'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0810'
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.33'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.044.11'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.630.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.849'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '49.389.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0946'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '270.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.92'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.137'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.50'
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.92'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.89'
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '128.61'
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.03'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0324'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.15'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.057.90'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.21'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.10'
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.90'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '59.10'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.20'
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) column (E)
$ws.Range("E2").Value = '  -0.76%  '
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("E6").Value = '  -0.79%  '
$ws.Range("E7").Value = '  -1.05%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -2.54%  '
$ws.Range("E10").Value = '  -3.20%  '
$ws.Range("E11").Value = '  -2.96%  '
$ws.Range("E12").Value = '  -1.24%  '
$ws.Range("E13").Value = '  +1.42%  '
$ws.Range("E14").Value = '  +0.00%  '
$ws.Range("E15").Value = '  -0.73%  '
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("E17").Value = '  -3.60%  '
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("E19").Value = '  -3.05%  '
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("E21").Value = '  -1.89%  '
$ws.Range("E22").Value = '  -1.50%  '
$ws.Range("E23").Value = '  -3.64%  '
$ws.Range("E24").Value = '  -5.61%  '
$ws.Range("E25").Value = '  -2.13%  '
$ws.Range("E26").Value = '  -2.69%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("E28").Value = '  +3.29%  '
$ws.Range("E29").Value = '  -1.34%  '
$ws.Range("E30").Value = '  -4.74%  '
$ws.Range("E31").Value = '  -5.66%  '
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("E33").Value = '  +1.07%  '
$ws.Range("E34").Value = '  +2.06%  '
$ws.Range("E35").Value = '  -0.23%  '
$ws.Range("E36").Value = '  -3.23%  '
$ws.Range("E37").Value = '  +2.64%  '
$ws.Range("E38").Value = '  -0.93%  '
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("E40").Value = '  +1.40%  '
$ws.Range("E41").Value = '  -1.57%  '
$ws.Range("E42").Value = '  -1.11%  '
$ws.Range("E43").Value = '  +3.21%  '
$ws.Range("E44").Value = '  -4.29%  '
$ws.Range("E45").Value = '  -0.54%  '
$ws.Range("E46").Value = '  -5.20%  '
$ws.Range("E47").Value = '  +6.16%  '
$ws.Range("E48").Value = '  -5.45%  '
$ws.Range("E49").Value = '  -2.06%  '
$ws.Range("E50").Value = '  +2.24%  '
$ws.Range("E51").Value = '  -3.90%  '
